# Adds two new columns, I (I0) and J (IF), to the sheet:
#  - I1/J1 get the header labels "I0"/"IF" (styled like the other headers)
#  - I2:J77 get per-row numeric values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header cells --------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, borders, centered) from the existing "IP"
# header cell (H1) onto the two new header cells without touching their
# values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Data values -----------------------------------------------------
$iValues = @(9,9,6,8,7,7,8,8,9,8,8,7,7,9,7,8,7,8,6,8,8,6,6,7,9,8,8,7,8,8,7,8,7,7,8,7,9,6,7,8,7,8,8,7,6,9,7,7,9,9,7,8,9,10,9,7,8,7,7,7,7,7,7,6,7,7,7,5,9,8,7,6,8,6,5,7)
$jValues = @(9,9,6,8,7,7,8,8,9,8,8,7,8,9,7,8,7,8,6,8,8,6,7,7,9,8,8,8,8,8,7,8,7,7,8,8,9,7,7,8,8,8,8,8,6,9,8,7,9,9,8,8,9,11,9,8,8,7,7,7,7,7,7,6,7,7,7,5,9,8,7,6,8,6,5,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
